$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 'He told the head-man the following story: The Sky-god has sent me; he says, I must take this single grain of corn and go with it somewhere, and when I come (anywhere), I must let this grain of corn sleep with the fowls. '
$ws.Range("B10").Value = 'The head-man said, " E! when fowls and corn sleep together, will it be well? " '
$ws.Range("B11").Value = 'Ananse replied, "All you have to do is to allow them to sleep together." '
$ws.Range("B12").Value = 'The head-man said, " Go and place it there." '
$ws.Range("B13").Value = 'Next morning early Ananse said he was going to look for it, but the (single grain of) corn was no longer there. '
$ws.Range("B14").Value = 'The Spider said, " Head-man, the fowls have eaten the Sky-god''s corn." '
$ws.Range("B15").Value = 'The head-man said, " I am not able to have any dispute about this, so take the fowl and go." '
$ws.Range("B16").Value = 'Ananse went and took one large fowl, and set out with it for another village. '
$ws.Range("B17").Value = 'When he reached there, he told the head-man a story, saying, "The Sky-god sent me, and I shall sleep here, and this fowl which I hold must sleep with the sheep." '
$ws.Range("B18").Value = 'The head-man said, "Bring it and put it there." '
$ws.Range("B19").Value = 'In the night, Ananse went to the sheep kraal and killed the fowl. '
$ws.Range("B20").Value = 'He took the fowl''s intestines and stuck them on the horns of a big ram. '
$ws.Range("B21").Value = 'Next morning early, he went to the head-man there and said, " Give me the fowl." '
$ws.Range("B22").Value = 'When they would have gone and fetched it for him, behold, the fowl was dead. '
$ws.Range("B23").Value = 'Ananse said, "Ah! that is the Sky-god''s fowl, which the sheep has killed." '
$ws.Range("B24").Value = 'The head-man said, "I am not able to have any dispute about it, so take a sheep and go." '
$ws.Range("B25").Value = 'Ananse took the sheep and-set out. '
$ws.Range("B26").Value = 'He reached another village. '
$ws.Range("B27").Value = 'He went to the head-man''s house.'
$ws.Range("B28").Value = 'The head-man said, "What news? " '
$ws.Range("B29").Value = 'The Spider replied, " The Sky-god has sent me, and I hold the sheep; it and the cattle are to sleep together." '
$ws.Range("B30").Value = 'The head-man said, " Ho! let it rest with the sheep." '
$ws.Range("B31").Value = 'The Spider said, " No, for this is the special sheep of the Sky-god, it rests with the cattle." '
$ws.Range("B32").Value = 'The head-man said, " Go and put it there, then." '
$ws.Range("B33").Value = 'The Spider went and put it there. '
$ws.Range("B34").Value = 'At night, Ananse went and killed the sheep, and put its blood on the head of one of the cattle.'
$ws.Range("B35").Value = 'Next morning early, Ananse said, " Head-man, give me my sheep." '
$ws.Range("B36").Value = 'They went to the kraal; there was the sheep, dead. '
$ws.Range("B37").Value = 'Ananse said, " Head-man, it is the Sky god''s sheep which the cattle have killed." '
$ws.Range("B38").Value = 'The head-man said, " I am not able to have any dispute with you about this matter, so take away the cow that killed it." '
$ws.Range("B39").Value = 'Ananse took the cow, and set off, and came to another village. '
$ws.Range("B40").Value = 'He took the cow and tied it up there, and went to sleep.'

# Rows whose text got shorter after the edit no longer need the taller
# wrapped-text row height, so autofit them back down to the default.
$ws.Rows.Item(10).AutoFit()
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(40).AutoFit()

# Restore the view state left by the author after the edit: scrolled down
# the sheet, selecting B42, with the zoom reset to 100%.
$ws.Range("B42").Select() | Out-Null
$excel.ActiveWindow.Zoom = 100
